# "Add files via upload" -- refresh the Football Manager team import list:
#   * fix the "Hudderfield" typo -> "Huddersfield" (keeps its spot among the
#     already-imported Premier League clubs)
#   * append the remaining Championship / lower-league clubs so the sheet
#     now runs from row 1 to row 36

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$teams = @(
    "Man Utd",
    "Man City",
    "Liverpool",
    "Leicester",
    "Chelsea",
    "Huddersfield",
    "Aston Villa",
    "Brighton",
    "Norwich",
    "Newcastle",
    "Sheff Utd",
    "West Ham",
    "Arsenal",
    "Leeds",
    "Nottingham Forest",
    "Burnley",
    "Middlesbrough",
    "Watford",
    "Bolton",
    "Sunderland",
    "West Bromwich",
    "Bournemouth",
    "Portsmouth",
    "Charlton",
    "Wigan",
    "Derby",
    "Ipswich",
    "Coventry",
    "Sheff Wed",
    "Swindon",
    "Barnsley",
    "Wolverhampton",
    "Milwall",
    "Oldham",
    "Luton",
    "Gillingham"
)

for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $teams[$i]
}

$last = $teams.Length
$ws.Range("A$last").Select() | Out-Null
